$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the stale "selection" element from the sheet view
$ws.Range("A1").Select()

# Remove the now-unused helper column G for every data row (kept on row 1)
$ws.Range("G2:G127").Clear()

# F2:F127 switch from the literal "1" placeholder text to the real stock count (numeric, centered)
$rng = $ws.Range("F2:F127")
$rng.NumberFormat = "#,##0"
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4160

$ws.Range("F2").Value = 15
$ws.Range("F3").Value = 33
$ws.Range("F4").Value = 10
$ws.Range("F5").Value = 34
$ws.Range("F6").Value = 13
$ws.Range("F7").Value = 5
$ws.Range("F8").Value = 227
$ws.Range("F9").Value = 88
$ws.Range("F10").Value = 95
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 43
$ws.Range("F13").Value = 58
$ws.Range("F14").Value = 54
$ws.Range("F15").Value = 12
$ws.Range("F16").Value = 5
$ws.Range("F17").Value = 2.8
$ws.Range("F18").Value = 4
$ws.Range("F19").Value = 5
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = 468
$ws.Range("F23").Value = 315
$ws.Range("F24").Value = 90
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("F28").Value = 3
$ws.Range("F29").Value = 1.9999999999999876
$ws.Range("F30").Value = 3.0000000000000133
$ws.Range("F31").Value = [double]"-1.2434497875801753e-14"
$ws.Range("F32").Value = 5.000000000000021
$ws.Range("F33").Value = 250
$ws.Range("F34").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("F36").Value = 158
$ws.Range("F37").Value = 17
$ws.Range("F38").Value = 46
$ws.Range("F39").Value = 9
$ws.Range("F40").Value = 5
$ws.Range("F41").Value = 33
$ws.Range("F42").Value = 0
$ws.Range("F43").Value = 16
$ws.Range("F44").Value = 18
$ws.Range("F45").Value = 6
$ws.Range("F46").Value = 3
$ws.Range("F47").Value = 35
$ws.Range("F48").Value = 33
$ws.Range("F49").Value = 9
$ws.Range("F50").Value = 43
$ws.Range("F51").Value = 0
$ws.Range("F52").Value = 0
$ws.Range("F53").Value = 9
$ws.Range("F54").Value = 22
$ws.Range("F55").Value = 18
$ws.Range("F56").Value = 0
$ws.Range("F57").Value = 0
$ws.Range("F58").Value = 6
$ws.Range("F59").Value = 17
$ws.Range("F60").Value = 8
$ws.Range("F61").Value = 3
$ws.Range("F62").Value = 13
$ws.Range("F63").Value = 11
$ws.Range("F64").Value = 0
$ws.Range("F65").Value = 2
$ws.Range("F66").Value = 4
$ws.Range("F67").Value = 16
$ws.Range("F68").Value = 21
$ws.Range("F69").Value = 12
$ws.Range("F70").Value = 5
$ws.Range("F71").Value = 10
$ws.Range("F72").Value = 13
$ws.Range("F73").Value = 12
$ws.Range("F74").Value = 80
$ws.Range("F75").Value = 14
$ws.Range("F76").Value = 17
$ws.Range("F77").Value = 11
$ws.Range("F78").Value = 13
$ws.Range("F79").Value = 0
$ws.Range("F80").Value = 75
$ws.Range("F81").Value = 45
$ws.Range("F82").Value = 9
$ws.Range("F83").Value = 14
$ws.Range("F84").Value = 33
$ws.Range("F85").Value = 24
$ws.Range("F86").Value = 27
$ws.Range("F87").Value = 30
$ws.Range("F88").Value = 38
$ws.Range("F89").Value = -6
$ws.Range("F90").Value = 76
$ws.Range("F91").Value = 58
$ws.Range("F92").Value = 52
$ws.Range("F93").Value = 2
$ws.Range("F94").Value = 0
$ws.Range("F95").Value = 88
$ws.Range("F96").Value = 84
$ws.Range("F97").Value = 25
$ws.Range("F98").Value = 0
$ws.Range("F99").Value = 88
$ws.Range("F100").Value = 8
$ws.Range("F101").Value = 0
$ws.Range("F102").Value = 2
$ws.Range("F103").Value = 52
$ws.Range("F104").Value = 9
$ws.Range("F105").Value = 53
$ws.Range("F106").Value = 26
$ws.Range("F107").Value = 3
$ws.Range("F108").Value = 0
$ws.Range("F109").Value = 93
$ws.Range("F110").Value = 67
$ws.Range("F111").Value = 0
$ws.Range("F112").Value = 11
$ws.Range("F113").Value = 9
$ws.Range("F114").Value = 8
$ws.Range("F115").Value = 20
$ws.Range("F116").Value = 0
$ws.Range("F117").Value = 12
$ws.Range("F118").Value = 15
$ws.Range("F119").Value = 2
$ws.Range("F120").Value = 24
$ws.Range("F121").Value = 8
$ws.Range("F122").Value = 68
$ws.Range("F123").Value = 24
$ws.Range("F124").Value = 25
$ws.Range("F125").Value = 130
$ws.Range("F126").Value = 51
$ws.Range("F127").Value = 27

# Rows 128/129 keep their text style but refresh the stock values
$ws.Range("F128").Value = "0"
$ws.Range("F129").Value = "15"

# Page setup entry added by the print-ready export
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
